# Insert a new weekly price record at row 349 on the single worksheet.
# This pushes the previously-existing rows 349-389 down to 350-390
# (keeping all their data/styles intact) and the newly opened row 349
# is populated with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("349:349").Insert()

$ws.Range("A349").Value = 4
$ws.Range("B349").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C349").Value = "Los Lagos"
$ws.Range("D349").Value = 45142
$ws.Range("E349").Value = 10
$ws.Range("F349").Value = 100112039
$ws.Range("G349").Value = "Ciboulette"
$ws.Range("H349").Value = "Sin especificar"
$ws.Range("I349").Value = "Primera"
$ws.Range("J349").Value = 240
$ws.Range("K349").Value = 3000
$ws.Range("L349").Value = 3500
$ws.Range("M349").Value = 3250
$ws.Range("N349").Value = "`$/docena de atados"
$ws.Range("O349").Value = "Región Metropolitana"
$ws.Range("P349").Value = 1083
$ws.Range("Q349").Value = 3
$ws.Range("R349").Value = "Hortaliza"
